# This script reproduces the "Add files via upload" commit for
# Business Analysis/data/user_4.xlsx:
#   - the vague "Afternoon (12pm - 6pm)" period value is replaced by the
#     shorter "afternoon" label
#   - the song/artist information that used to live only in the combined
#     "song-artist" column is additionally split into two brand new
#     columns, J ("song") and K ("artist"), for the header row and the
#     one data row that is populated
#   - the selection / column sizing bookkeeping that Excel keeps in the
#     worksheet is refreshed to reflect the newly added column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the "period" value in the single populated data row.
$ws.Range("G2").Value = "afternoon"

# Add the two new columns with their headers ...
$ws.Range("J1").Value = "song"
$ws.Range("K1").Value = "artist"

# ... and their values for the existing data row, split out of the
# "Calma-Pedro Capo" song-artist value already present in column H.
$ws.Range("J2").Value = "Calma"
$ws.Range("K2").Value = "Pedro Capo"

# Column H ("song-artist") now gets an explicit best-fit-style width,
# matching the author's saved worksheet.
$ws.Columns.Item(8).ColumnWidth = 15.8333333333333

# The author's last selection in the sheet was on the newly added K7 cell.
$ws.Range("K7").Select() | Out-Null
